$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamps = @(
    "Fri Sep 29 11:29:55 EDT 2023",
    "Fri Sep 29 11:30:08 EDT 2023",
    "Fri Sep 29 11:30:21 EDT 2023",
    "Fri Sep 29 11:30:34 EDT 2023",
    "Fri Sep 29 11:30:47 EDT 2023",
    "Fri Sep 29 11:31:00 EDT 2023",
    "Fri Sep 29 11:31:13 EDT 2023",
    "Fri Sep 29 11:31:26 EDT 2023",
    "Fri Sep 29 11:31:39 EDT 2023",
    "Fri Sep 29 11:31:52 EDT 2023",
    "Fri Sep 29 11:32:05 EDT 2023",
    "Fri Sep 29 11:32:18 EDT 2023",
    "Fri Sep 29 11:32:31 EDT 2023",
    "Fri Sep 29 11:32:45 EDT 2023",
    "Fri Sep 29 11:32:58 EDT 2023",
    "Fri Sep 29 11:33:11 EDT 2023",
    "Fri Sep 29 11:33:25 EDT 2023",
    "Fri Sep 29 11:33:38 EDT 2023",
    "Fri Sep 29 11:33:52 EDT 2023",
    "Fri Sep 29 11:34:05 EDT 2023",
    "Fri Sep 29 11:34:18 EDT 2023",
    "Fri Sep 29 11:34:31 EDT 2023",
    "Fri Sep 29 11:34:44 EDT 2023",
    "Fri Sep 29 11:34:57 EDT 2023",
    "Fri Sep 29 11:35:10 EDT 2023",
    "Fri Sep 29 11:35:24 EDT 2023",
    "Fri Sep 29 11:35:37 EDT 2023",
    "Fri Sep 29 11:35:50 EDT 2023",
    "Fri Sep 29 11:36:04 EDT 2023",
    "Fri Sep 29 11:36:17 EDT 2023",
    "Fri Sep 29 11:36:31 EDT 2023",
    "Fri Sep 29 11:36:44 EDT 2023"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
}
